$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text in C1
$ws.Range("C1").Value = "DisplayName"

# Fill in the previously-empty row 2 with new data
$ws.Range("A2").Value = "Running"
$ws.Range("B2").Value = "AarSvc_427f0fb"
$ws.Range("C2").Value = "Agent Activation Runtime_427f0fb"

# Correct the mangled (mojibake) registered-trademark text in C3
$ws.Range("C3").Value = "IntelÂ® SGX AESM"

# Append two new rows after the existing data
$ws.Range("A4").Value = "Stopped"
$ws.Range("B4").Value = "AJRouter"
$ws.Range("C4").Value = "AllJoyn Router Service"

$ws.Range("A5").Value = "Stopped"
$ws.Range("B5").Value = "ALG"
$ws.Range("C5").Value = "Application Layer Gateway Service"
